$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update row 1 values (existing columns B/C change, new columns D/E added)
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "10.0.0.0/24"
$ws.Range("C1").Value = "10.0.0.34/32"

# New header cells D1 / E1 - copy formatting from an existing header cell (A1)
# then set their value, so the style (bold/centered) matches the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("D1").Value = "10.0.0.49/32"
$ws.Range("E1").Value = "10.0.0.23/32"

# ---------------------------------------------------------------------------
# 2. Remove the old row 3 ("*" / "" / "Protocol : * - Port : *") - its content
#    effectively gets folded into row 2 below.
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# 3. Rebuild row 2 content
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "*"
$ws.Range("B2").Value = "Protocol : * - Port : *"

# New data cells C2 / D2 / E2 - copy formatting from B2 (style used for the
# protocol/port cells) before filling them in with their content.
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)

$ws.Range("C2").Value = "Protocol : TCP - Port : 80`nProtocol : TCP - Port : 443"
$ws.Range("D2").Value = "Protocol : TCP - Port : 143`nProtocol : TCP - Port : 993"
$ws.Range("E2").Value = "Protocol : UDP - Port : 389"

# Row 2 grows taller to fit the two-line protocol/port entries.
$ws.Rows.Item(2).RowHeight = 30

# ---------------------------------------------------------------------------
# 4. Resize columns: B shrinks, C keeps the "wide" width, D/E (new) match C.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 22.1666666
$ws.Columns.Item(3).ColumnWidth = 26.1666666
$ws.Columns.Item(4).ColumnWidth = 26.1666666
$ws.Columns.Item(5).ColumnWidth = 26.1666666
